$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the three "ID / Due Date" rows (A3:B5) and the three
# "Invoice Number / Invoice Date / Company Name / Total Due" rows
# (C6:F8) into three complete rows (3-5) with refreshed data, then drop
# the now-empty rows 6-8.
#
# All of these values must stay plain TEXT (they are shared strings in
# the source file, not numbers/dates), so every cell is temporarily
# switched to the "@" text format before the value is assigned and then
# restored to the default/general format so no stray cell style sticks
# around.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = ""
}

Set-TextValue $ws.Range("A3") "bzfw5m88g0fymg6lk5cle"
Set-TextValue $ws.Range("B3") "16-07-2024"
Set-TextValue $ws.Range("C3") "284213"
Set-TextValue $ws.Range("D3") "2019-06-03"
Set-TextValue $ws.Range("E3") "Aenean LLC"
Set-TextValue $ws.Range("F3") "9778.40"

Set-TextValue $ws.Range("A4") "agkf70jn9satt1rxtvyy6"
Set-TextValue $ws.Range("B4") "07-07-2024"
Set-TextValue $ws.Range("C4") "284221"
Set-TextValue $ws.Range("D4") "2019-06-20"
Set-TextValue $ws.Range("E4") "Aenean LLC"
Set-TextValue $ws.Range("F4") "6300.00"

Set-TextValue $ws.Range("A5") "wj34k48z92mgkik0lpdt1g"
Set-TextValue $ws.Range("B5") "25-07-2024"
Set-TextValue $ws.Range("C5") "284232"
Set-TextValue $ws.Range("D5") "2019-06-15"
Set-TextValue $ws.Range("E5") "Aenean LLC"
Set-TextValue $ws.Range("F5") "1009.80"

# The old data used to live in C6:F8 - it has now been folded into rows
# 3-5 above, so clear it (rows 6-8 become empty and drop out of the
# sheet's used range).
$ws.Range("C6:F8").Clear()

# Update the saved selection/active cell.
$ws.Range("B4").Select()

$wb.Save()
